$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.55
$ws.Range("I2").Value = 2.75
$ws.Range("AB3").Value = 17
$ws.Range("E4").Value = "Den Bosch"
$ws.Range("F4").Value = "Telstar"
$ws.Range("G4").Value = 2.8
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 2.45
$ws.Range("J4").Value = 1.05
$ws.Range("K4").Value = 11
$ws.Range("L4").Value = 1.29
$ws.Range("M4").Value = 3.5
$ws.Range("N4").Value = 1.9
$ws.Range("O4").Value = 1.9
$ws.Range("P4").Value = 1.36
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = 1.7
$ws.Range("S4").Value = 2.05
$ws.Range("T4").Value = 9.5
$ws.Range("U4").Value = 15
$ws.Range("V4").Value = 11
$ws.Range("W4").Value = 29
$ws.Range("X4").Value = 21
$ws.Range("Y4").Value = 29
$ws.Range("Z4").Value = 11
$ws.Range("AA4").Value = 6.5
$ws.Range("AB4").Value = 13
$ws.Range("AC4").Value = 41
$ws.Range("AD4").Value = 201
$ws.Range("AE4").Value = 8.5
$ws.Range("AF4").Value = 12
$ws.Range("AG4").Value = 9.5
$ws.Range("AH4").Value = 23
$ws.Range("AI4").Value = 19
$ws.Range("AJ4").Value = 26
$ws.Range("G5").Value = 1.5
$ws.Range("I5").Value = 8
$ws.Range("W5").Value = 9.5
$ws.Range("AB5").Value = 26
$ws.Range("AJ5").Value = 81
$ws.Range("G8").Value = 7
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 1.5
$ws.Range("W8").Value = 67
$ws.Range("J9").Value = 1.06
$ws.Range("K9").Value = 10
$ws.Range("G13").Value = 1.62
$ws.Range("L13").Value = 1.17
$ws.Range("M13").Value = 5
$ws.Range("N13").Value = 1.57
$ws.Range("O13").Value = 2.35
$ws.Range("P13").Value = 1.25
$ws.Range("Q13").Value = 3.75
$ws.Range("T13").Value = 9.5
$ws.Range("U13").Value = 9.5
$ws.Range("AA13").Value = 8.5
$ws.Range("G18").Value = 1.42
$ws.Range("H18").Value = 4.05
$ws.Range("I18").Value = 6.2
$ws.Range("N18").Value = 1.65
$ws.Range("O18").Value = 1.98
$ws.Range("T18").Value = 6
$ws.Range("U18").Value = 5.9
$ws.Range("V18").Value = 6.9
$ws.Range("W18").Value = 8
$ws.Range("X18").Value = 9.5
$ws.Range("Y18").Value = 20
$ws.Range("Z18").Value = 11.75
$ws.Range("AA18").Value = 7.1
$ws.Range("AB18").Value = 14.5
$ws.Range("AC18").Value = 60
$ws.Range("AD18").Value = 400
$ws.Range("AE18").Value = 14.5
$ws.Range("AF18").Value = 32
$ws.Range("AG18").Value = 16
$ws.Range("AH18").Value = 100
$ws.Range("AI18").Value = 50
$ws.Range("AJ18").Value = 45
$ws.Range("G19").Value = 1.1
$ws.Range("H19").Value = 7.3
$ws.Range("I19").Value = 15
$ws.Range("T19").Value = 9.25
$ws.Range("U19").Value = 6
$ws.Range("V19").Value = 10
$ws.Range("W19").Value = 5.6
$ws.Range("X19").Value = 9.25
$ws.Range("Y19").Value = 28
$ws.Range("Z19").Value = 22
$ws.Range("AA19").Value = 15.5
$ws.Range("AB19").Value = 30
$ws.Range("AC19").Value = 120
$ws.Range("AE19").Value = 40
$ws.Range("AF19").Value = 120
$ws.Range("AG19").Value = 45
$ws.Range("AH19").Value = 500
$ws.Range("AI19").Value = 200
$ws.Range("AJ19").Value = 120
$ws.Range("K20").Value = 10
$ws.Range("L20").Value = 1.29
$ws.Range("M20").Value = 3.5
$ws.Range("N20").Value = 1.98
$ws.Range("O20").Value = 1.83
$ws.Range("Z20").Value = 10
$ws.Range("AC20").Value = 41
$ws.Range("AJ20").Value = 29
$ws.Range("J24").Value = 1.03
$ws.Range("K24").Value = 10.5
$ws.Range("J25").Value = 1.04
$ws.Range("K25").Value = 9
$ws.Range("J26").Value = 1.02
$ws.Range("K26").Value = 12
$ws.Range("H29").Value = 3.5
$ws.Range("I29").Value = 3.75
$ws.Range("R29").Value = 1.7
$ws.Range("T29").Value = 7.4
$ws.Range("V29").Value = 8.75
$ws.Range("AA29").Value = 7.1
$ws.Range("AC29").Value = 65
$ws.Range("AE29").Value = 11.25
$ws.Range("AF29").Value = 23
$ws.Range("AH29").Value = 60
